$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a row of C:G time values (copying the number format used
# by the existing time columns) and an optional Notes (H) string.
# NOTE: positional parameters only -- named parameter binding (-Foo bar)
# is not reliable in this host, so every call below passes args in order.
function Set-CortisolRow($Row, $AtAwake, $Min45, $Hours8, $Hours12, $BeforeBed, $Note) {
    $ws.Range("C$Row").Value = $AtAwake
    $ws.Range("D$Row").Value = $Min45
    $ws.Range("E$Row").Value = $Hours8
    $ws.Range("F$Row").Value = $Hours12
    $ws.Range("G$Row").Value = $BeforeBed
    $ws.Range("C$Row`:G$Row").NumberFormat = $ws.Range("C2").NumberFormat

    if ($Note) {
        $ws.Range("H$Row").Value = $Note
    }
}

# A handful of the source values are tiny (times just after midnight) and
# only round-trip exactly in scientific notation -- this PS host's parser
# doesn't accept bare `1.23E-4` numeric literals, so build those few via a
# string->double cast stashed in a variable first.
$v_F7  = [double]"6.9444444444444447E-4"
$v_G7  = [double]"6.9444444444444441E-3"
$v_G10 = [double]"6.25E-2"
$v_G22 = [double]"3.1944444444444442E-2"
$v_G33 = [double]"7.4305555555555555E-2"
$v_G39 = [double]"8.3333333333333332E-3"
$v_G42 = [double]"3.0555555555555555E-2"

# --- Week 3 data for participants that were missing it ---

# Row 4: Redwood (participant 2), week 3
Set-CortisolRow 4 0.29930555555555555 0.33333333333333331 0.64583333333333337 0.80902777777777779 0.98958333333333337

# Row 7: Granite (participant 3), week 3
Set-CortisolRow 7 0.46250000000000002 0.50069444444444444 0.83402777777777781 $v_F7 $v_G7 "(final 2 times are technically the next day, and only 10 minutes apart)"

# Row 10: Glacier (participant 4), week 3
Set-CortisolRow 10 0.50694444444444442 0.53819444444444442 0.875 0 $v_G10 "(final 2 times are technically the next day)"

# Row 19: Marigold (participant 7), week 3
Set-CortisolRow 19 0.39583333333333331 0.42708333333333331 0.72916666666666663 0.89583333333333337 0.95833333333333337

# Row 22: Solstice (participant 8), week 3
Set-CortisolRow 22 0.31597222222222221 0.3611111111111111 0.69097222222222221 0.81597222222222221 $v_G22 "(final time is technically the next day)"

# Row 26: Cascade (participant 10), week 1 -- note text corrected to the
# more specific "only 15 minutes apart" wording.
$ws.Range("H26").Value = "(final 2 times are technically the next day, and only 15 minutes apart)"

# Row 33: Eclipse (participant 12), week 2
Set-CortisolRow 33 0.2326388888888889 0.2638888888888889 0.56597222222222221 0.73263888888888884 $v_G33 "(final time is technically the next day)"

# Row 36: Quartz (participant 13), week 2
Set-CortisolRow 36 0.41666666666666669 0.44791666666666669 0.71875 0.91249999999999998 0.11458333333333333 "(final time is technically the next day)"

# Row 39: Tundra (participant 14), week 2
Set-CortisolRow 39 0.3347222222222222 0.3659722222222222 0.70763888888888893 0.83680555555555558 $v_G39 "(final time is technically the next day)"

# Row 42: Harbor (participant 15), week 2
Set-CortisolRow 42 0.2986111111111111 0.3298611111111111 0.63888888888888884 0.85 $v_G42 "(final time is technically the next day)"

# Row 45: Obsidian (participant 16), week 2
Set-CortisolRow 45 0.32847222222222222 0.34375 0.68055555555555558 0.875 0.91666666666666663

# --- View state: scroll down a bit and select H43 ---
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H43").Select()
